$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.03819759486660548
$ws.Range("H2").Value = -9.65277585525263
$ws.Range("I2").Value = 33.00337386403519
$ws.Range("G3").Value = 0.0595672719658755
$ws.Range("H3").Value = 22.31384977114901
$ws.Range("G4").Value = 0.007016730674758977
$ws.Range("H4").Value = 268.545987671781
$ws.Range("G5").Value = -0.00558174888139492
$ws.Range("H5").Value = -247.0293182823429
$ws.Range("G6").Value = 0.06261000079036701
$ws.Range("H6").Value = 80.59059778421118
$ws.Range("G7").Value = 0.08855204987330606
$ws.Range("H7").Value = 66.48581291063182
$ws.Range("G8").Value = 0.01490658233141812
$ws.Range("H8").Value = 179.2028081301959
$ws.Range("G9").Value = 0.001690300959692852
$ws.Range("H9").Value = 107.8516741809202
$ws.Range("G10").Value = -0.04729306482862562
$ws.Range("H10").Value = 34.94885456082012
$ws.Range("G11").Value = -0.07418717711622526
$ws.Range("H11").Value = 19.37144077220124
$ws.Range("G12").Value = -0.2340960532931304
$ws.Range("H12").Value = 4.237428573188897
$ws.Range("G13").Value = -0.3111000085225923
$ws.Range("H13").Value = -13.20416371354739
$ws.Range("G14").Value = -0.05075162001597489
$ws.Range("H14").Value = -36.80716808910736
$ws.Range("G15").Value = -0.003085842890363873
$ws.Range("H15").Value = 91.1253314452846
$ws.Range("G16").Value = 0.1346668713025791
$ws.Range("H16").Value = 7.468533416354393
$ws.Range("G17").Value = 0.1475871522120627
$ws.Range("H17").Value = 5.231112538889964
$ws.Range("G18").Value = 0.1331620833676958
$ws.Range("H18").Value = 6.757937365241736
$ws.Range("G19").Value = 0.1278542482069725
$ws.Range("H19").Value = -4.023301599360891
$ws.Range("G20").Value = 0.03678928601984719
$ws.Range("H20").Value = 7.144586591373874
$ws.Range("G21").Value = 0.05267653309215381
$ws.Range("H21").Value = -9.240768809702251
$ws.Range("G22").Value = -0.05266148268420764
$ws.Range("H22").Value = 34.04782441385697
$ws.Range("G23").Value = -0.07016960362575986
$ws.Range("H23").Value = -12.18801564343637
$ws.Range("G24").Value = 0.1145094226431902
$ws.Range("H24").Value = -3.052196883576929
$ws.Range("G25").Value = 0.1232162706591249
$ws.Range("H25").Value = -2.341375812972784
$ws.Range("G26").Value = 0.04823283847361042
$ws.Range("H26").Value = -2.960056588517513
$ws.Range("G27").Value = 0.09273467425940869
$ws.Range("H27").Value = 6.989321260092233
$ws.Range("G28").Value = -0.07609160330851496
$ws.Range("H28").Value = -19.65210499696565
$ws.Range("G29").Value = -0.09550856950380406
$ws.Range("H29").Value = -34.19360530611367
$ws.Range("G30").Value = 0.07855707453166022
$ws.Range("H30").Value = 23.30777639333505
$ws.Range("G31").Value = 0.05835684704198968
$ws.Range("H31").Value = -3.670513852526665
$ws.Range("G32").Value = 0.1014777440641883
$ws.Range("H32").Value = 3.270536841167464
$ws.Range("G33").Value = 0.06338803230345763
$ws.Range("H33").Value = -22.96498820713023
$ws.Range("G34").Value = -0.008802242539078907
$ws.Range("H34").Value = -133.7830342715564
$ws.Range("G35").Value = 0.007474504490262423
$ws.Range("H35").Value = 166.676466436311
$ws.Range("G36").Value = 0.00476771868138504
$ws.Range("H36").Value = 785.852090427739
$ws.Range("G37").Value = -0.0006052085631946093
$ws.Range("H37").Value = 95.17925927146115
$ws.Range("G38").Value = 0.1020358789660517
$ws.Range("H38").Value = -4.868097603937708
$ws.Range("G39").Value = 0.1022818874773153
$ws.Range("H39").Value = 19.40061381057111
$ws.Range("G40").Value = 0.004039995908217165
$ws.Range("H40").Value = 36.01501543374192
$ws.Range("G41").Value = 0.03635336617999158
$ws.Range("H41").Value = 142.4083539905599
$ws.Range("G42").Value = 0.09575626016116683
$ws.Range("H42").Value = -5.130885148000504
$ws.Range("G43").Value = 0.1211154539897585
$ws.Range("H43").Value = 0.8081524711119243
$ws.Range("G44").Value = 0.03007479626620252
$ws.Range("H44").Value = -15.72757159719789
$ws.Range("G45").Value = 0.0332216816693883
$ws.Range("H45").Value = 102.9403067904029
$ws.Range("G46").Value = 0.05049762348090088
$ws.Range("H46").Value = 39.34792225167033
$ws.Range("G47").Value = 0.07233176916606916
$ws.Range("H47").Value = 43.40071189418956
$ws.Range("G48").Value = 0.03644658757496372
$ws.Range("H48").Value = -14.80720389621818
$ws.Range("G49").Value = 0.0709187596797867
$ws.Range("H49").Value = 2.07648240889519
$ws.Range("G50").Value = 0.01319781775301627
$ws.Range("H50").Value = -23.59148948388412
$ws.Range("G51").Value = 0.01502316724646558
$ws.Range("H51").Value = -22.83826292731912
$ws.Range("G52").Value = -0.1006740260478277
$ws.Range("H52").Value = 2.749629870665762
$ws.Range("G53").Value = -0.07752343921450217
$ws.Range("H53").Value = 16.05930129476169
$ws.Range("G54").Value = 0.06769015725334654
$ws.Range("H54").Value = -7.429302636817595
$ws.Range("G55").Value = 0.09616158660878361
$ws.Range("H55").Value = 55.22036172530475
$ws.Range("G56").Value = 0.02780428365732081
$ws.Range("H56").Value = -20.53470274091964
$ws.Range("G57").Value = 0.0282656451653466
$ws.Range("H57").Value = 389.5740604584028
$ws.Range("G58").Value = 0.03863652560886717
$ws.Range("H58").Value = 54.48055973573867
$ws.Range("G59").Value = 0.02757144665196607
$ws.Range("H59").Value = 16.4402314035841
$ws.Range("G60").Value = 0.0249463967996456
$ws.Range("H60").Value = -23.1061869173646
$ws.Range("G61").Value = 0.02511361585224851
$ws.Range("H61").Value = 98.40124350948713
$ws.Range("G62").Value = 0.05398668124619469
$ws.Range("H62").Value = -10.56272389239195
$ws.Range("G63").Value = 0.03678473407739007
$ws.Range("H63").Value = 12.8729854382194
$ws.Range("G64").Value = 0.03614603284979163
$ws.Range("H64").Value = -10.80819374318273
$ws.Range("G65").Value = 0.05966828493075703
$ws.Range("H65").Value = 6.432822329222314
$ws.Range("G66").Value = 0.09762585163865245
$ws.Range("H66").Value = 4.352291187526736
$ws.Range("G67").Value = 0.1149906491444582
$ws.Range("H67").Value = -0.3949939282026267
$ws.Range("G68").Value = -0.03012942153526091
$ws.Range("H68").Value = 13.54661670416778
$ws.Range("G69").Value = -0.01799198648811699
$ws.Range("H69").Value = 15.21953053959942
$ws.Range("G70").Value = 0.07652342148165765
$ws.Range("H70").Value = -17.39366078780789
$ws.Range("G71").Value = 0.09795915617870746
$ws.Range("H71").Value = 7.400996564895028
$ws.Range("G72").Value = -0.05013110604563263
$ws.Range("H72").Value = 10.60950646764424
$ws.Range("G73").Value = -0.06288130032358182
$ws.Range("H73").Value = 14.7517726370817
$ws.Range("G74").Value = 0.1033392411128273
$ws.Range("H74").Value = 3.394239822676357
$ws.Range("G75").Value = 0.117859872218662
$ws.Range("H75").Value = 21.00103154261331
$ws.Range("G76").Value = 0.005878292642772768
$ws.Range("H76").Value = -77.01151705826661
$ws.Range("G77").Value = 0.02920149829515547
$ws.Range("H77").Value = 106.9701268717602
$ws.Range("G78").Value = 0.1005745059376496
$ws.Range("H78").Value = 56.47051829243934
$ws.Range("G79").Value = 0.06603887814941581
$ws.Range("H79").Value = -13.91536424727973
$ws.Range("G80").Value = -0.1522676894068221
$ws.Range("H80").Value = 8.055577336363996
$ws.Range("G81").Value = -0.1568173966277395
$ws.Range("H81").Value = 25.35775479544268
$ws.Range("G82").Value = 0.1187731201093882
$ws.Range("H82").Value = 3.555325365377098
$ws.Range("G83").Value = 0.1944623116636404
$ws.Range("H83").Value = 9.25940796262074
$ws.Range("G84").Value = 0.08157184049789443
$ws.Range("H84").Value = 242.1975059148194
$ws.Range("G85").Value = 0.06021927539860066
$ws.Range("H85").Value = -2.203238477907317
